$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header: "Lembaga" -> "Divisi/Unit"
$ws.Range("G1").Value = "Divisi/Unit"

# Data rows: reassign each member's "Lembaga A"/"Lembaga B" value to a specific "Divisi N"
$divisiValues = @{
    2  = "Divisi 1"
    3  = "Divisi 2"
    4  = "Divisi 3"
    5  = "Divisi 4"
    6  = "Divisi 5"
    7  = "Divisi 6"
    8  = "Divisi 7"
    9  = "Divisi 8"
    10 = "Divisi 9"
    11 = "Divisi 2"
    12 = "Divisi 3"
    13 = "Divisi 4"
    14 = "Divisi 5"
    15 = "Divisi 6"
    16 = "Divisi 7"
    17 = "Divisi 8"
    18 = "Divisi 9"
    19 = "Divisi 10"
    20 = "Divisi 11"
    21 = "Divisi 12"
}

foreach ($row in $divisiValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $divisiValues[$row]
}

# Move the active selection to match the saved view state
$null = $ws.Range("G24").Select()
